$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price & 1h volume change columns),
# matching the refreshed data pulled from coinranking.com.
# Columns D (Price) and E (Volume 1h) are kept as plain text
# (mirrors the original workbook, which stores these as text too)
# so numeric-looking strings such as "217.97" are not coerced to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.127.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.656.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.11%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5162"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2573"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.65%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06438"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.87"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07761"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.654.59"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.297"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.883.31"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5530"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8039"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.32"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.171.02"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.52"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.06"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.917"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.88"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1160"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.972"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.73"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05293"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.253"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.359"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.233"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.571"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.761"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9223"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5686"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.161.72"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +10.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01589"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8372"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.646"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.83"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.793.78"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4510"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.94"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.008"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.899"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.90%  "
